$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rsquo = [char]0x2019

$a2 = "('Faerie Conclave', ['Land', 'Faerie Conclave enters the battlefield tapped.', '{T}: Add {U}.', '{1}{U}: Faerie Conclave becomes a 2/1 blue Faerie creature with flying until end of turn. It" + $rsquo + "s still a land.'])"
$a3 = "('Treetop Village', ['Land', 'Treetop Village enters the battlefield tapped.', '{T}: Add {G}.', '{1}{G}: Treetop Village becomes a 3/3 green Ape creature with trample until end of turn. It" + $rsquo + "s still a land. (It can deal excess combat damage to the player or planeswalker it" + $rsquo + "s attacking.)'])"

$ws.Range("A2").Value = $a2
$ws.Range("A3").Value = $a3

$ws.Range("A4:A11").EntireRow.Delete()
